$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Step 1: extend invalidCredentialTest (sheet1) with two more data rows ---
$ws1.Range("A3").Value = "bala"
$ws1.Range("B3").Value = "bala123"
$ws1.Range("C3").Value = "Dutch"
$ws1.Range("D3").Value = "Invalid username or password"

$ws1.Range("A4").Value = "john123"
$ws1.Range("B4").Value = "john123"
$ws1.Range("C4").Value = "Danish"
$ws1.Range("D4").Value = "Invalid username or password"

# --- Step 2: create checkLanguageTest sheet right after sheet1 and fill its
#     first 5 data rows (matches shared-string creation order from the
#     original authoring session). NOTE: worksheet object variables here are
#     position-bound, not identity-bound - once another sheet is inserted
#     before/at this position the old variable silently starts resolving to
#     whatever sheet now sits at that index. So every later touch of this
#     sheet re-fetches it by name instead of reusing this variable. ---
$wsLang = $wb.Worksheets.Add($null, $ws1)
$wsLang.Name = "checkLanguageTest"

$wsLang.Range("A1").Value = "ExpectedLanguage"
$wsLang.Range("A2").Value = "Albanian"
$wsLang.Range("A3").Value = "Amharic"
$wsLang.Range("A4").Value = "Arabic"
$wsLang.Range("A5").Value = "Armenian"

# --- Step 3: create checkVersionNumberTest sheet right after sheet1 (i.e.
#     pushing checkLanguageTest one tab to the right), matching the final
#     target tab order: invalidCredentialTest, checkVersionNumberTest,
#     checkLanguageTest, addPatientTest. ---
$wsVer = $wb.Worksheets.Add($null, $ws1)
$wsVer.Name = "checkVersionNumberTest"

$wsVer.Range("A1").Value = "Username"
$wsVer.Range("B1").Value = "Password"
$wsVer.Range("C1").Value = "Language"
$wsVer.Range("D1").Value = "Expected Version"

$wsVer.Range("A2").Value = "admin"
$wsVer.Range("B2").Value = "pass"
$wsVer.Range("C2").Value = "English (Indian)"

$wsVer.Range("A3").Value = "accountant"
$wsVer.Range("B3").Value = "accountant"
$wsVer.Range("C3").Value = "English (Indian)"

$wsVer.Range("D2").Value = "Version Number: v6.0.0 (2)"
$wsVer.Range("D3").Value = "Version Number: v6.0.0 (2)"

# --- Step 4: go back to checkLanguageTest (re-fetched by name - see note in
#     step 2) and append the last language row. ---
$wsLang = $wb.Worksheets.Item("checkLanguageTest")
$wsLang.Range("A6").Value = "Japanese"

# --- Step 5: create addPatientTest sheet at the end ---
$wsPatient = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsPatient.Name = "addPatientTest"

# DOB column (F) is stored as text (numFmtId 49 "@") in the target workbook, so
# set that number format before writing into it - otherwise a date-looking
# string like "2021-12-16" gets auto-converted to a date serial number.
$wsPatient.Range("F1:F2").NumberFormat = "@"

$wsPatient.Range("A1").Value = "Username"
$wsPatient.Range("B1").Value = "Password"
$wsPatient.Range("C1").Value = "Language"
$wsPatient.Range("D1").Value = "First Nme"
$wsPatient.Range("E1").Value = "Last Name"
$wsPatient.Range("F1").Value = "DOB"
$wsPatient.Range("G1").Value = "Gender"
$wsPatient.Range("H1").Value = "Expected Alert"
$wsPatient.Range("I1").Value = "Expected Patient Detail"

$wsPatient.Range("A2").Value = "admin"
$wsPatient.Range("B2").Value = "pass"
$wsPatient.Range("C2").Value = "English (Indian)"
$wsPatient.Range("D2").Value = "John"
$wsPatient.Range("E2").Value = "Wick"
$wsPatient.Range("F2").Value = "2021-12-16"
$wsPatient.Range("G2").Value = "Male"
$wsPatient.Range("H2").Value = "Tobacco"
$wsPatient.Range("I2").Value = "Medical Record Dashboard - John Wick"

# --- Step 6: selections / active sheet / view state (re-fetch every sheet by
#     name so nothing relies on a possibly-stale position-bound variable) ---
$wb.Worksheets.Item("invalidCredentialTest").Range("A1:C1").Select()
$wb.Worksheets.Item("checkVersionNumberTest").Range("A2:C2").Select()
$wb.Worksheets.Item("checkLanguageTest").Range("A7").Select()
$wb.Worksheets.Item("addPatientTest").Range("I10").Select()

$wb.Worksheets.Item("addPatientTest").Activate()
